$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: AD1 "Wins", AE1 "Losses", AF1 "Ties" ---
# Copy the formatting of the existing last header cell (AC1 - bold,
# centered/top-aligned, thin boxed border) onto the three new header
# cells so they visually match the rest of row 1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Season record values for every player row (2-50) ---
# The team finished the season 89-73-0; stamp that record onto every
# player row in the new Wins / Losses / Ties columns.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 89   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 73   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-50"
